$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = 45961
$ws.Range("B19").Value = 5594
$ws.Range("C19").Value = 3925
$ws.Range("D19").Value = 3630
$ws.Range("E19").Value = 222
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = 5
$ws.Range("I19").Value = 1

$ws.Range("A19:I19").Select()
